# Slide 5, "Content Placeholder 2" shape: the paragraph
#   "Advance " + "User Interface "
# (two separate runs) should become a single run:
#   "Advance User Interface "
# keeping the second run's formatting (lang="en-US" dirty="0" smtClean="0").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the two runs making up the last bullet of the placeholder text:
# "Advance " (8 chars) starting at 79, followed by "User Interface " (15 chars) at 87.
$runAdvance = $tr.Characters(79, 8)
$runUI = $tr.Characters(87, 15)

# Sanity-check we grabbed the right runs before mutating anything.
if ($runAdvance.Text -ne "Advance " -or $runUI.Text -ne "User Interface ") {
    throw "Unexpected text layout: [$($runAdvance.Text)] / [$($runUI.Text)]"
}

# Clear the first run entirely so it (and its run properties) is dropped,
# leaving the second run - whose rPr carries dirty="0" - as the sole run.
$runAdvance.Text = ""

# Re-acquire the (now shifted) range covering the merged text and set the
# full combined string on it, preserving the surviving run's formatting.
$tr2 = $sh.TextFrame.TextRange
$merged = $tr2.Characters(79, 15)
$merged.Text = "Advance User Interface "
